# "Last Day at R Systems" - PO_Detail sheet gets 4 more PO rows appended
# below the existing row 2, so the reconciler has several POs ("Stuck on
# UNTPRG-Confirm PO") to work through instead of just one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO_Detail")

# Insert 4 fresh rows under row 2 (keeps row 2's look, which we'll
# re-point below) so the table grows from 1 data row to 5 data rows.
$ws.Rows("3:6").Insert() | Out-Null

# The insert cloned row 2's formatting into columns A and D for every
# new row; A and D are only needed on specific rows below, so drop the
# inherited formatting from the rows that don't use them.
$ws.Range("A3:A6").Clear() | Out-Null
$ws.Range("D3:D5").Clear() | Out-Null

# Row 2: unit number cleared (no longer tied to a single VIN/unit), PO
# is now "01000996", comment explains it's stuck on PO confirmation.
$ws.Range("A2").ClearContents() | Out-Null
$ws.Range("B2").Value = "'1"
$ws.Range("C2").Value = "'01000996"
$ws.Range("D2").Value = "Stuck on UNTPRG-Confirm PO"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null

# Rows 3-6: additional PO numbers (group "2") awaiting reconciliation.
$ws.Range("B3").Value = "'2"
$ws.Range("C3").Value = "'01000949"
$ws.Range("B3").Style = "Normal"

$ws.Range("B4").Value = "'2"
$ws.Range("C4").Value = "'01000942"
$ws.Range("B4").Style = "Normal"

$ws.Range("B5").Value = "'2"
$ws.Range("C5").Value = "'01000943"
$ws.Range("B5").Style = "Normal"

$ws.Range("B6").Value = "'2"
$ws.Range("C6").Value = "'01001011"
$ws.Range("D6").Value = "Stuck on UNTPRG-Confirm PO"
$ws.Range("B6").Style = "Normal"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null

# Leave the selection on the newest PO entry, like the author did.
$ws.Range("C3").Select() | Out-Null
